$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 37500
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 37500
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 37500
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -38472

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3352.2075
$ws.Range("I76").Value = 2959.7073
$ws.Range("J76").Value = 4693.25
$ws.Range("K76").Value = 2959.7073
$ws.Range("L76").Value = 4693.25
$ws.Range("M76").Value = -2644.7073
$ws.Range("N76").Value = -5323.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3352.2075
$ws.Range("I79").Value = 2959.7073
$ws.Range("J79").Value = 4693.25
$ws.Range("K79").Value = 2959.7073
$ws.Range("L79").Value = 4693.25
$ws.Range("M79").Value = -1867.7073
$ws.Range("N79").Value = -6877.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 8446.682000000001
$ws.Range("I113").Value = 2673.5715
$ws.Range("J113").Value = 11140.8
$ws.Range("K113").Value = 2673.5715
$ws.Range("L113").Value = 11140.8
$ws.Range("M113").Value = 580.4285
$ws.Range("N113").Value = -17648.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1325.711
$ws.Range("I137").Value = 937.20514
$ws.Range("K137").Value = 2811.61542
$ws.Range("M137").Value = -261.6154200000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3348.5847
$ws.Range("I138").Value = 1657.8636
$ws.Range("J138").Value = 4213.6045
$ws.Range("K138").Value = 4973.5908
$ws.Range("L138").Value = 12640.8135
$ws.Range("M138").Value = 166.4092000000001
$ws.Range("N138").Value = -22920.8135

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1047.6842
$ws.Range("I61").Value = 895.5
$ws.Range("K61").Value = 895.5
$ws.Range("M61").Value = -683.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1141.4286
$ws.Range("J74").Value = 1145.6
$ws.Range("L74").Value = 1145.6
$ws.Range("N74").Value = -2893.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1141.4286
$ws.Range("J77").Value = 1145.6
$ws.Range("L77").Value = 5728
$ws.Range("N77").Value = -14464

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1914.1945
$ws.Range("I132").Value = 1252.3334
$ws.Range("J132").Value = 2386.9524
$ws.Range("K132").Value = 3757.0002
$ws.Range("L132").Value = 7160.8572
$ws.Range("M132").Value = -1227.0002
$ws.Range("N132").Value = -12220.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1047.6842
$ws.Range("I136").Value = 895.5
$ws.Range("K136").Value = 2686.5
$ws.Range("M136").Value = -136.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 49985
$ws.Range("J139").Value = 49985
$ws.Range("L139").Value = 49985
$ws.Range("N139").Value = -60265

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 45000
$ws.Range("J138").Value = 45000
$ws.Range("L138").Value = 45000
$ws.Range("N138").Value = -55280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50471

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 49999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2141849.8
$ws.Range("I99").Value = 3201874.8
$ws.Range("J99").Value = 21800
$ws.Range("K99").Value = 3201874.8
$ws.Range("L99").Value = 21800
$ws.Range("M99").Value = -3200376.8
$ws.Range("N99").Value = -24796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2141849.8
$ws.Range("I126").Value = 3201874.8
$ws.Range("J126").Value = 21800
$ws.Range("K126").Value = 9605624.399999999
$ws.Range("L126").Value = 65400
$ws.Range("M126").Value = -9603154.399999999
$ws.Range("N126").Value = -70340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1098000
$ws.Range("I5").Value = 990
$ws.Range("K5").Value = 2970
$ws.Range("M5").Value = -2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 2378.5715
$ws.Range("I99").Value = 750
$ws.Range("J99").Value = 3600
$ws.Range("K99").Value = 2250
$ws.Range("L99").Value = 10800
$ws.Range("M99").Value = -4
$ws.Range("N99").Value = -15292

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 42583.832
$ws.Range("I121").Value = 125315
$ws.Range("J121").Value = 1218.25
$ws.Range("K121").Value = 375945
$ws.Range("L121").Value = 3654.75
$ws.Range("M121").Value = -374635
$ws.Range("N121").Value = -6274.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1200.5
$ws.Range("I122").Value = 371
$ws.Range("J122").Value = 2237.375
$ws.Range("K122").Value = 3339
$ws.Range("L122").Value = 20136.375
$ws.Range("M122").Value = -889
$ws.Range("N122").Value = -25036.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10477.272
$ws.Range("I131").Value = 383.33334
$ws.Range("J131").Value = 14262.5
$ws.Range("K131").Value = 1150.00002
$ws.Range("L131").Value = 42787.5
$ws.Range("M131").Value = 3889.99998
$ws.Range("N131").Value = -52867.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 2136.875
$ws.Range("I133").Value = 2365
$ws.Range("K133").Value = 7095
$ws.Range("M133").Value = -2035

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1098000
$ws.Range("I135").Value = 990
$ws.Range("K135").Value = 8910
$ws.Range("M135").Value = -6375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 66704.25
$ws.Range("J141").Value = 66704.25
$ws.Range("L141").Value = 66704.25
$ws.Range("N141").Value = -77064.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 169668.17
$ws.Range("I7").Value = 251500
$ws.Range("J7").Value = 6004.5
$ws.Range("K7").Value = 251500
$ws.Range("L7").Value = 6004.5
$ws.Range("M7").Value = -251388
$ws.Range("N7").Value = -6228.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 14000
$ws.Range("J62").Value = 14000
$ws.Range("L62").Value = 14000
$ws.Range("N62").Value = -15248

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 14000
$ws.Range("J65").Value = 14000
$ws.Range("L65").Value = 42000
$ws.Range("N65").Value = -48240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 169668.17
$ws.Range("I126").Value = 251500
$ws.Range("J126").Value = 6004.5
$ws.Range("K126").Value = 754500
$ws.Range("L126").Value = 18013.5
$ws.Range("M126").Value = -752030
$ws.Range("N126").Value = -22953.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 44182.9
$ws.Range("J138").Value = 44182.9
$ws.Range("L138").Value = 44182.9
$ws.Range("N138").Value = -54462.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 66550.5
$ws.Range("J140").Value = 66550.5
$ws.Range("L140").Value = 66550.5
$ws.Range("N140").Value = -76910.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 57142.223
$ws.Range("J141").Value = 57142.223
$ws.Range("L141").Value = 57142.223
$ws.Range("N141").Value = -67502.223

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1703.3269
$ws.Range("I132").Value = 1682.738
$ws.Range("J132").Value = 1789.8
$ws.Range("K132").Value = 5048.214
$ws.Range("L132").Value = 5369.4
$ws.Range("M132").Value = -2518.214
$ws.Range("N132").Value = -10429.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1428.3158
$ws.Range("I136").Value = 712.63635
$ws.Range("J136").Value = 2412.375
$ws.Range("K136").Value = 2137.90905
$ws.Range("L136").Value = 7237.125
$ws.Range("M136").Value = 412.0909499999998
$ws.Range("N136").Value = -12337.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 50715
$ws.Range("J141").Value = 50715
$ws.Range("L141").Value = 50715
$ws.Range("N141").Value = -61075
